$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the existing row 1206, pushing all
# subsequent rows (old 1206-1272) down to 1210-1276.
$ws.Rows("1206:1209").Insert()

# Row 1206 - new market entry (Morada(o) / Primera, Peru)
$ws.Cells.Item(1206, 1).Value = 3
$ws.Cells.Item(1206, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1206, 3).Value = "Coquimbo"
$ws.Cells.Item(1206, 4).Value = 44706
$ws.Cells.Item(1206, 5).Value = 5
$ws.Cells.Item(1206, 6).Value = 100112004
$ws.Cells.Item(1206, 7).Value = "Cebolla"
$ws.Cells.Item(1206, 8).Value = "Morada(o)"
$ws.Cells.Item(1206, 9).Value = "Primera"
$ws.Cells.Item(1206, 10).Value = 105
$ws.Cells.Item(1206, 11).Value = 11500
$ws.Cells.Item(1206, 12).Value = 12000
$ws.Cells.Item(1206, 13).Value = 11762
$ws.Cells.Item(1206, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(1206, 15).Value = "Perú"
$ws.Cells.Item(1206, 16).Value = 653
$ws.Cells.Item(1206, 17).Value = 18
$ws.Cells.Item(1206, 18).Value = "Hortaliza"

# Row 1207 - new market entry (Morada(o) / Segunda, Peru)
$ws.Cells.Item(1207, 1).Value = 3
$ws.Cells.Item(1207, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1207, 3).Value = "Coquimbo"
$ws.Cells.Item(1207, 4).Value = 44706
$ws.Cells.Item(1207, 5).Value = 5
$ws.Cells.Item(1207, 6).Value = 100112004
$ws.Cells.Item(1207, 7).Value = "Cebolla"
$ws.Cells.Item(1207, 8).Value = "Morada(o)"
$ws.Cells.Item(1207, 9).Value = "Segunda"
$ws.Cells.Item(1207, 10).Value = 55
$ws.Cells.Item(1207, 11).Value = 9000
$ws.Cells.Item(1207, 12).Value = 9000
$ws.Cells.Item(1207, 13).Value = 9000
$ws.Cells.Item(1207, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(1207, 15).Value = "Perú"
$ws.Cells.Item(1207, 16).Value = 500
$ws.Cells.Item(1207, 17).Value = 18
$ws.Cells.Item(1207, 18).Value = "Hortaliza"

# Row 1208 - new market entry (Sin especificar / 1a (guarda), Provincia de Quillota)
$ws.Cells.Item(1208, 1).Value = 3
$ws.Cells.Item(1208, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1208, 3).Value = "Coquimbo"
$ws.Cells.Item(1208, 4).Value = 44706
$ws.Cells.Item(1208, 5).Value = 5
$ws.Cells.Item(1208, 6).Value = 100112004
$ws.Cells.Item(1208, 7).Value = "Cebolla"
$ws.Cells.Item(1208, 8).Value = "Sin especificar"
$ws.Cells.Item(1208, 9).Value = "1a (guarda)"
$ws.Cells.Item(1208, 10).Value = 185
$ws.Cells.Item(1208, 11).Value = 6000
$ws.Cells.Item(1208, 12).Value = 6500
$ws.Cells.Item(1208, 13).Value = 6257
$ws.Cells.Item(1208, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(1208, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1208, 16).Value = 348
$ws.Cells.Item(1208, 17).Value = 18
$ws.Cells.Item(1208, 18).Value = "Hortaliza"

# Row 1209 - new market entry (Sin especificar / 2a (guarda), Provincia de Quillota)
$ws.Cells.Item(1209, 1).Value = 3
$ws.Cells.Item(1209, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1209, 3).Value = "Coquimbo"
$ws.Cells.Item(1209, 4).Value = 44706
$ws.Cells.Item(1209, 5).Value = 5
$ws.Cells.Item(1209, 6).Value = 100112004
$ws.Cells.Item(1209, 7).Value = "Cebolla"
$ws.Cells.Item(1209, 8).Value = "Sin especificar"
$ws.Cells.Item(1209, 9).Value = "2a (guarda)"
$ws.Cells.Item(1209, 10).Value = 95
$ws.Cells.Item(1209, 11).Value = 5000
$ws.Cells.Item(1209, 12).Value = 5000
$ws.Cells.Item(1209, 13).Value = 5000
$ws.Cells.Item(1209, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(1209, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1209, 16).Value = 278
$ws.Cells.Item(1209, 17).Value = 18
$ws.Cells.Item(1209, 18).Value = "Hortaliza"
